# Auto-generated edit script: updates crypto price/volume table values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.439.81"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "3.308.52"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'185.78"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'576.52"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").Value = "'0.410"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "3.887.26"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'27.44"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "67.631.45"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "3.331.35"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "'443.81"
$ws.Range("E18").Value = "  +6.61%  "
$ws.Range("D19").Value = "'5.68"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'7.75"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "'74.07"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").Value = "3.455.80"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").Value = "'9.05"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "'22.93"
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").Value = "'162.77"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'27.23"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "'4.47"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "2.752.32"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").Value = "'6.26"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'24.84"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").Value = "'40.19"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0671"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'326.91"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").Value = "'0.991"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  -1.40%  "
